# Applies the Typhon_Profits.xlsx profit-recalculation update described by the
# commit "chore: update Sheets via scheduled runner": a batch of refreshed
# market-price-driven values (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 185.88235
$ws.Range("I33").Value = 191.25
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 191.25
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = 37.75
$ws.Range("N33").Value = -558

# Row 127
$ws.Range("H127").Value = 1095.7391
$ws.Range("I127").Value = 561.25
$ws.Range("J127").Value = 1380.8
$ws.Range("K127").Value = 1683.75
$ws.Range("L127").Value = 4142.4
$ws.Range("M127").Value = 3276.25
$ws.Range("N127").Value = -14062.4

# Row 129
$ws.Range("H129").Value = 182868.38
$ws.Range("J129").Value = 189759.7
$ws.Range("L129").Value = 569279.1000000001
$ws.Range("N129").Value = -579279.1000000001

# Row 132
$ws.Range("H132").Value = 2158.76
$ws.Range("I132").Value = 2450.8
$ws.Range("J132").Value = 990.6
$ws.Range("K132").Value = 7352.400000000001
$ws.Range("L132").Value = 2971.8
$ws.Range("M132").Value = -4822.400000000001
$ws.Range("N132").Value = -8031.8

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 122.25
$ws.Range("I5").Value = 144.5
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 144.5
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -32.5
$ws.Range("N5").Value = -324

# Row 45
$ws.Range("H45").Value = 2550.1738
$ws.Range("I45").Value = 1932.8
$ws.Range("K45").Value = 1932.8
$ws.Range("M45").Value = -1555.8

# Row 61
$ws.Range("H61").Value = 1352.6666
$ws.Range("I61").Value = 1350.05
$ws.Range("J61").Value = 1365.75
$ws.Range("K61").Value = 1350.05
$ws.Range("L61").Value = 1365.75
$ws.Range("M61").Value = -1138.05
$ws.Range("N61").Value = -1789.75

# Row 110
$ws.Range("H110").Value = 591.0769
$ws.Range("I110").Value = 520.6667
$ws.Range("K110").Value = 520.6667
$ws.Range("M110").Value = 1524.3333

# Row 122
$ws.Range("H122").Value = 1221.5
$ws.Range("I122").Value = 1067.6522
$ws.Range("J122").Value = 2106.125
$ws.Range("K122").Value = 3202.9566
$ws.Range("L122").Value = 6318.375
$ws.Range("M122").Value = -752.9566
$ws.Range("N122").Value = -11218.375

# Row 132
$ws.Range("H132").Value = 22406.72
$ws.Range("I132").Value = 2132.8
$ws.Range("J132").Value = 103502.4
$ws.Range("K132").Value = 6398.400000000001
$ws.Range("L132").Value = 310507.2
$ws.Range("M132").Value = -3868.400000000001
$ws.Range("N132").Value = -315567.2

# Row 136
$ws.Range("H136").Value = 1352.6666
$ws.Range("I136").Value = 1350.05
$ws.Range("J136").Value = 1365.75
$ws.Range("K136").Value = 4050.15
$ws.Range("L136").Value = 4097.25
$ws.Range("M136").Value = -1500.15
$ws.Range("N136").Value = -9197.25

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 122.25
$ws.Range("I4").Value = 144.5
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 144.5
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -29.5
$ws.Range("N4").Value = -330

# Row 80
$ws.Range("H80").Value = 1162.3103
$ws.Range("I80").Value = 1544.2727
$ws.Range("J80").Value = 928.8889
$ws.Range("K80").Value = 1544.2727
$ws.Range("L80").Value = 928.8889
$ws.Range("M80").Value = -546.2727
$ws.Range("N80").Value = -2924.8889

# Row 83
$ws.Range("H83").Value = 1162.3103
$ws.Range("I83").Value = 1544.2727
$ws.Range("J83").Value = 928.8889
$ws.Range("K83").Value = 7721.363499999999
$ws.Range("L83").Value = 4644.444500000001
$ws.Range("M83").Value = -2729.363499999999
$ws.Range("N83").Value = -14628.4445

# Row 94
$ws.Range("H94").Value = 931.46155
$ws.Range("I94").Value = 723.2222
$ws.Range("J94").Value = 1400
$ws.Range("K94").Value = 723.2222
$ws.Range("L94").Value = 1400
$ws.Range("M94").Value = -272.2222
$ws.Range("N94").Value = -2302

# Row 107
$ws.Range("H107").Value = 578.05884
$ws.Range("I107").Value = 609.2143
$ws.Range("J107").Value = 432.66666
$ws.Range("K107").Value = 609.2143
$ws.Range("L107").Value = 432.66666
$ws.Range("M107").Value = 1310.7857
$ws.Range("N107").Value = -4272.66666

# Row 134
$ws.Range("H134").Value = 2522.8447
$ws.Range("I134").Value = 2435.6482
$ws.Range("J134").Value = 3700
$ws.Range("K134").Value = 7306.944600000001
$ws.Range("L134").Value = 11100
$ws.Range("M134").Value = -4771.944600000001
$ws.Range("N134").Value = -16170

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3381.1904
$ws.Range("I31").Value = 3446.818
$ws.Range("J31").Value = 3357.9033
$ws.Range("K31").Value = 3446.818
$ws.Range("L31").Value = 3357.9033
$ws.Range("M31").Value = -3151.818
$ws.Range("N31").Value = -3947.9033

# Row 34
$ws.Range("H34").Value = 3381.1904
$ws.Range("I34").Value = 3446.818
$ws.Range("J34").Value = 3357.9033
$ws.Range("K34").Value = 3446.818
$ws.Range("L34").Value = 3357.9033
$ws.Range("M34").Value = -3244.818
$ws.Range("N34").Value = -3761.9033

# Row 58
$ws.Range("H58").Value = 17717.645
$ws.Range("I58").Value = 1562.2667
$ws.Range("J58").Value = 32863.312
$ws.Range("K58").Value = 1562.2667
$ws.Range("L58").Value = 32863.312
$ws.Range("M58").Value = -1359.2667
$ws.Range("N58").Value = -33269.312

# Row 99
$ws.Range("H99").Value = 4000
$ws.Range("I99").Value = 3083.3333
$ws.Range("J99").Value = 6750
$ws.Range("K99").Value = 3083.3333
$ws.Range("L99").Value = 6750
$ws.Range("M99").Value = -1585.3333
$ws.Range("N99").Value = -9746

# Row 107
$ws.Range("H107").Value = 1376.4
$ws.Range("I107").Value = 510.33334
$ws.Range("J107").Value = 2085
$ws.Range("K107").Value = 510.33334
$ws.Range("L107").Value = 2085
$ws.Range("M107").Value = 1409.66666
$ws.Range("N107").Value = -5925

# Row 126
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 3083.3333
$ws.Range("J126").Value = 6750
$ws.Range("K126").Value = 9249.999899999999
$ws.Range("L126").Value = 20250
$ws.Range("M126").Value = -6779.999899999999
$ws.Range("N126").Value = -25190

# Row 132
$ws.Range("H132").Value = 2468.3823
$ws.Range("I132").Value = 1753.88
$ws.Range("J132").Value = 4453.1113
$ws.Range("K132").Value = 5261.64
$ws.Range("L132").Value = 13359.3339
$ws.Range("M132").Value = -2731.64
$ws.Range("N132").Value = -18419.3339

# Row 134
$ws.Range("H134").Value = 1059.8334
$ws.Range("I134").Value = 946.9032
$ws.Range("J134").Value = 1760
$ws.Range("K134").Value = 2840.7096
$ws.Range("L134").Value = 5280
$ws.Range("M134").Value = -305.7096000000001
$ws.Range("N134").Value = -10350

# Row 136
$ws.Range("H136").Value = 17717.645
$ws.Range("I136").Value = 1562.2667
$ws.Range("J136").Value = 32863.312
$ws.Range("K136").Value = 4686.800099999999
$ws.Range("L136").Value = 98589.93599999999
$ws.Range("M136").Value = -2136.800099999999
$ws.Range("N136").Value = -103689.936

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = -2888
$ws.Range("N4").Value = -3224

# Row 80
$ws.Range("H80").Value = 2559.4
$ws.Range("I80").Value = 1666
$ws.Range("J80").Value = 3899.5
$ws.Range("K80").Value = 4998
$ws.Range("L80").Value = 11698.5
$ws.Range("M80").Value = -4062
$ws.Range("N80").Value = -13570.5

# Row 83
$ws.Range("H83").Value = 2559.4
$ws.Range("I83").Value = 1666
$ws.Range("J83").Value = 3899.5
$ws.Range("K83").Value = 14994
$ws.Range("L83").Value = 35095.5
$ws.Range("M83").Value = -10314
$ws.Range("N83").Value = -44455.5

# Row 113
$ws.Range("H113").Value = 733.1875
$ws.Range("I113").Value = 593.2857
$ws.Range("K113").Value = 1779.8571
$ws.Range("M113").Value = 390.1428999999998

# Row 117
$ws.Range("H117").Value = 1035.75
$ws.Range("J117").Value = 814.25
$ws.Range("L117").Value = 2442.75
$ws.Range("N117").Value = -9326.75

# Row 121
$ws.Range("H121").Value = 1466
$ws.Range("J121").Value = 1544.7273
$ws.Range("L121").Value = 4634.1819
$ws.Range("N121").Value = -7254.1819

# Row 131
$ws.Range("H131").Value = 763.61
$ws.Range("J131").Value = 759.28864
$ws.Range("L131").Value = 2277.86592
$ws.Range("N131").Value = -12357.86592

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1246.3182
$ws.Range("I97").Value = 1173.2778
$ws.Range("J97").Value = 1575
$ws.Range("K97").Value = 1173.2778
$ws.Range("L97").Value = 1575
$ws.Range("M97").Value = -677.2778000000001
$ws.Range("N97").Value = -2567

# Row 123
$ws.Range("H123").Value = 9281.75
$ws.Range("I123").Value = 3820
$ws.Range("J123").Value = 13183
$ws.Range("K123").Value = 3820
$ws.Range("L123").Value = 13183
$ws.Range("M123").Value = -1370
$ws.Range("N123").Value = -18083

# Row 132
$ws.Range("H132").Value = 30385.223
$ws.Range("I132").Value = 2409.3333
$ws.Range("J132").Value = 86337
$ws.Range("K132").Value = 7227.999899999999
$ws.Range("L132").Value = 259011
$ws.Range("M132").Value = -4697.999899999999
$ws.Range("N132").Value = -264071

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1059.4445
$ws.Range("I22").Value = 1564.5
$ws.Range("J22").Value = 655.4
$ws.Range("K22").Value = 1564.5
$ws.Range("L22").Value = 655.4
$ws.Range("M22").Value = -1269.5
$ws.Range("N22").Value = -1245.4

# Row 27
$ws.Range("H27").Value = 1059.4445
$ws.Range("I27").Value = 1564.5
$ws.Range("J27").Value = 655.4
$ws.Range("K27").Value = 1564.5
$ws.Range("L27").Value = 655.4
$ws.Range("M27").Value = -1457.5
$ws.Range("N27").Value = -869.4

# Row 74
$ws.Range("H74").Value = 36199.5
$ws.Range("I74").Value = 25197
$ws.Range("J74").Value = 38400
$ws.Range("K74").Value = 25197
$ws.Range("L74").Value = 38400
$ws.Range("M74").Value = -24199
$ws.Range("N74").Value = -40396

# Row 77
$ws.Range("H77").Value = 36199.5
$ws.Range("I77").Value = 25197
$ws.Range("J77").Value = 38400
$ws.Range("K77").Value = 75591
$ws.Range("L77").Value = 115200
$ws.Range("M77").Value = -70599
$ws.Range("N77").Value = -125184

# Row 93
$ws.Range("H93").Value = 3490
$ws.Range("I93").Value = 3490
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3490
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2242
$ws.Range("N93").ClearContents()

# Row 132
$ws.Range("H132").Value = 710907.6
$ws.Range("I132").Value = 1205955.2
$ws.Range("K132").Value = 3617865.6
$ws.Range("M132").Value = -3615335.6

# Row 136
$ws.Range("H136").Value = 1316.0588
$ws.Range("I136").Value = 1316.0588
$ws.Range("K136").Value = 3948.1764
$ws.Range("M136").Value = -1398.1764

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1265.9459
$ws.Range("I132").Value = 797.72
$ws.Range("J132").Value = 2241.4167
$ws.Range("K132").Value = 2393.16
$ws.Range("L132").Value = 6724.250100000001
$ws.Range("M132").Value = 136.8400000000001
$ws.Range("N132").Value = -11784.2501

# Row 136
$ws.Range("H136").Value = 21068712
$ws.Range("I136").Value = 28674840
$ws.Range("J136").Value = 5593.077
$ws.Range("K136").Value = 86024520
$ws.Range("L136").Value = 16779.231
$ws.Range("M136").Value = -86021970
$ws.Range("N136").Value = -21879.231
